$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 201
$ws1.Range("F4").Value = 83
$ws1.Range("F5").Value = 1671
$ws1.Range("F6").Value = 3266
$ws1.Range("F7").Value = 859
$ws1.Range("F8").Value = 2068
$ws1.Range("F9").Value = 1985
$ws1.Range("F10").Value = 1024
$ws1.Range("F11").Value = 358
$ws1.Range("F13").Value = 1618
$ws1.Range("F14").Value = 349
$ws1.Range("F18").Value = 84
$ws1.Range("F19").Value = 1454
$ws1.Range("F20").Value = 533
$ws1.Range("F22").Value = 326
$ws1.Range("F23").Value = 10790
$ws1.Range("F24").Value = 11784
$ws1.Range("F25").Value = 860
$ws1.Range("F26").Value = 666
$ws1.Range("F27").Value = 1850
$ws1.Range("F28").Value = 151

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 63

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 63
$ws4.Range("F4").Value = 201
$ws4.Range("F6").Value = 83
$ws4.Range("F7").Value = 1671
$ws4.Range("F8").Value = 3266
$ws4.Range("F9").Value = 859
$ws4.Range("F10").Value = 2068
$ws4.Range("F11").Value = 1985
$ws4.Range("F12").Value = 1024
$ws4.Range("F13").Value = 358
$ws4.Range("G14").Value = 128
$ws4.Range("F15").Value = 1618
$ws4.Range("F16").Value = 349
$ws4.Range("F22").Value = 84
$ws4.Range("F23").Value = 1454
$ws4.Range("F24").Value = 533
$ws4.Range("F26").Value = 326
$ws4.Range("F27").Value = 10791
$ws4.Range("F28").Value = 11784
$ws4.Range("F29").Value = 860
$ws4.Range("F30").Value = 666
$ws4.Range("F31").Value = 1850
$ws4.Range("F34").Value = 151
